$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1) Append the 6 new observation rows (2011-2016) for Congo, DRC GDP per
#    Capita. Columns A-D are plain numbers / reused text; column E (the
#    data value) is handled below together with the other years because it
#    must be written as *text* (it is stored as a shared string in the
#    source workbook, even though it looks numeric).
# ---------------------------------------------------------------------------
$newYears = @(2011, 2012, 2013, 2014, 2015, 2016)
for ($i = 0; $i -lt $newYears.Length; $i++) {
    $r = 63 + $i
    $ws.Range("A$r").Value = 180
    $ws.Range("B$r").Value = "Congo, DRC"
    $ws.Range("C$r").Value = "GDP per Capita"
    $ws.Range("D$r").Value = $newYears[$i]
}

# ---------------------------------------------------------------------------
# 2) Updated GDP per Capita series (1950-2016), row 2 .. row 68. Every
#    existing year's figure was revised and six new years were appended.
#    The values must land as text (shared-string) cells, matching the
#    source data's typing, so each one is staged as a formula in a scratch
#    column, copied, and pasted back as values - which preserves the text
#    type without leaving any NumberFormat/quote-prefix style behind.
# ---------------------------------------------------------------------------
$dataValues = @("909","995","1065","1105","1148","1170","1224","1237","1176","1191","1192","1041","1223","1248","1191","1176","1224","1176","1189","1262","1224","1267","1235","1305","1312","1213","1111","1086","993","961","953","950","920","909","932","909","923","918","894","861","813","737.934490092681","656.691651156199","563.71093100535","525.388717427869","531.990921495922","523.868062041872","489.642164758035","472.923964475332","452.004421422413","414.046947048094","403.6018922542","413.825212735492","436.028486254343","463.996975838844","490.690494007437","515.258365852912","546.036795512043","578.848721083029","595.016055121906","637.129849905218","681","710","751","802","836","836")

for ($i = 0; $i -lt $dataValues.Length; $i++) {
    $r = 2 + $i
    $ws.Range("ZZ$r").Formula = '="' + $dataValues[$i] + '"'
}

$lastRow = 1 + $dataValues.Length
$ws.Range("ZZ2:ZZ$lastRow").Copy()
$ws.Range("E2:E$lastRow").PasteSpecial(-4163)
$ws.Range("ZZ2:ZZ$lastRow").Clear()

Write-Host "Work Week and Social Spending data refreshed"
